# Update the cryptos price/volume table to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells are plain-text (inlineStr) in the source sheet; force text storage so
# numeric-looking values (e.g. "588.91", "0.250") are not reinterpreted as
# numbers, then restore the default "Normal" style so no stray number-format
# style is left behind on the cell.
$targetCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "E26",
    "D27",
    "E27",
    "E28",
    "D29",
    "E29",
    "D30",
    "E30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "E34",
    "E35",
    "D36",
    "E36",
    "D37",
    "E37",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "B47",
    "C47",
    "D47",
    "E47",
    "B48",
    "C48",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "E51"
)
foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value2 = "62.696.16"
$ws.Range("E2").Value2 = "  -1.69%  "
$ws.Range("D3").Value2 = "3.183.58"
$ws.Range("E3").Value2 = "  -3.70%  "
$ws.Range("D5").Value2 = "588.91"
$ws.Range("E5").Value2 = "  -2.35%  "
$ws.Range("D6").Value2 = "135.83"
$ws.Range("E6").Value2 = "  -4.56%  "
$ws.Range("D7").Value2 = "0.999"
$ws.Range("E7").Value2 = "  -0.15%  "
$ws.Range("D8").Value2 = "3.184.21"
$ws.Range("E8").Value2 = "  -3.61%  "
$ws.Range("D9").Value2 = "0.507"
$ws.Range("E9").Value2 = "  -2.52%  "
$ws.Range("D10").Value2 = "0.142"
$ws.Range("E10").Value2 = "  -4.67%  "
$ws.Range("D11").Value2 = "5.27"
$ws.Range("D12").Value2 = "0.454"
$ws.Range("E12").Value2 = "  -3.64%  "
$ws.Range("D13").Value2 = "0.0000236"
$ws.Range("E13").Value2 = "  -4.91%  "
$ws.Range("D14").Value2 = "33.35"
$ws.Range("E14").Value2 = "  -3.78%  "
$ws.Range("D15").Value2 = "3.699.95"
$ws.Range("E15").Value2 = "  -3.86%  "
$ws.Range("E16").Value2 = "  -1.52%  "
$ws.Range("D17").Value2 = "3.174.95"
$ws.Range("E17").Value2 = "  -4.04%  "
$ws.Range("D18").Value2 = "62.601.33"
$ws.Range("E18").Value2 = "  -1.97%  "
$ws.Range("D19").Value2 = "6.56"
$ws.Range("E19").Value2 = "  -4.74%  "
$ws.Range("D20").Value2 = "457.05"
$ws.Range("E20").Value2 = "  -4.82%  "
$ws.Range("D21").Value2 = "13.92"
$ws.Range("E21").Value2 = "  -1.70%  "
$ws.Range("D22").Value2 = "0.706"
$ws.Range("E22").Value2 = "  -3.75%  "
$ws.Range("D23").Value2 = "7.64"
$ws.Range("E23").Value2 = "  -4.56%  "
$ws.Range("D24").Value2 = "13.44"
$ws.Range("E24").Value2 = "  +0.15%  "
$ws.Range("D25").Value2 = "83.62"
$ws.Range("E25").Value2 = "  -1.66%  "
$ws.Range("E26").Value2 = "  +0.01%  "
$ws.Range("D27").Value2 = "2.69"
$ws.Range("E27").Value2 = "  -2.80%  "
$ws.Range("E28").Value2 = "  -0.05%  "
$ws.Range("D29").Value2 = "6.90"
$ws.Range("E29").Value2 = "  -5.19%  "
$ws.Range("D30").Value2 = "7.78"
$ws.Range("E30").Value2 = "  -4.18%  "
$ws.Range("D31").Value2 = "2.03"
$ws.Range("E31").Value2 = "  -6.41%  "
$ws.Range("D32").Value2 = "27.43"
$ws.Range("E32").Value2 = "  -6.60%  "
$ws.Range("D33").Value2 = "0.105"
$ws.Range("E33").Value2 = "  -0.73%  "
$ws.Range("E34").Value2 = "  -5.76%  "
$ws.Range("E35").Value2 = "  -5.73%  "
$ws.Range("D36").Value2 = "5.92"
$ws.Range("E36").Value2 = "  -0.84%  "
$ws.Range("D37").Value2 = "51.04"
$ws.Range("E37").Value2 = "  -3.48%  "
$ws.Range("D38").Value2 = "0.0₃0696"
$ws.Range("E38").Value2 = "  -7.00%  "
$ws.Range("D39").Value2 = "0.0386"
$ws.Range("E39").Value2 = "  -3.81%  "
$ws.Range("D40").Value2 = "408.16"
$ws.Range("E40").Value2 = "  -4.55%  "
$ws.Range("D41").Value2 = "2.70"
$ws.Range("E41").Value2 = "  -1.16%  "
$ws.Range("D42").Value2 = "2.870.70"
$ws.Range("E42").Value2 = "  -5.54%  "
$ws.Range("D43").Value2 = "8.02"
$ws.Range("E43").Value2 = "  -4.43%  "
$ws.Range("D44").Value2 = "0.112"
$ws.Range("E44").Value2 = "  +0.22%  "
$ws.Range("D45").Value2 = "36.56"
$ws.Range("E45").Value2 = "  +3.12%  "
$ws.Range("D46").Value2 = "0.250"
$ws.Range("E46").Value2 = "  -5.64%  "
$ws.Range("B47").Value2 = "Fetch.AI"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value2 = "2.14"
$ws.Range("E47").Value2 = "  -2.11%  "
$ws.Range("B48").Value2 = "USDe"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value2 = "0.999"
$ws.Range("E48").Value2 = "  -0.06%  "
$ws.Range("D49").Value2 = "124.99"
$ws.Range("E49").Value2 = "  +0.52%  "
$ws.Range("D50").Value2 = "25.50"
$ws.Range("E50").Value2 = "  -3.20%  "
$ws.Range("E51").Value2 = "  -3.33%  "

foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).Style = "Normal"
}
